# Apply updated cryptocurrency price/volume figures (and the three-row
# reshuffles among rows 43-45 and 50-51) per the Sun Feb 18 2024 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed by Excel as a number
# (e.g. "354.68"); force them to stay plain text via a leading apostrophe,
# then drop back to the default cell style so no explicit style id sticks.
$textForcedCells = @(
    'D5',
    'D6',
    'D10',
    'D12',
    'D13',
    'D14',
    'D20',
    'D21',
    'D23',
    'D24',
    'D29',
    'D32',
    'D38',
    'D40',
    'D43',
    'D44',
    'D45',
    'D46',
    'D48',
    'D50',
    'D51',
)

$ws.Range('D2').Value = '51.758.25'
$ws.Range('E2').Value = '  -0.01%  '
$ws.Range('D3').Value = '2.803.81'
$ws.Range('E3').Value = '  +0.74%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = "'354.68"
$ws.Range('E5').Value = '  -0.52%  '
$ws.Range('D6').Value = "'111.62"
$ws.Range('E6').Value = '  +2.22%  '
$ws.Range('E7').Value = '  +0.64%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('E9').Value = '  +7.80%  '
$ws.Range('D10').Value = "'40.36"
$ws.Range('E10').Value = '  +1.66%  '
$ws.Range('E11').Value = '  -1.46%  '
$ws.Range('D12').Value = "'0.0839"
$ws.Range('E12').Value = '  -0.91%  '
$ws.Range('D13').Value = "'20.01"
$ws.Range('E13').Value = '  +2.41%  '
$ws.Range('D14').Value = "'7.76"
$ws.Range('E14').Value = '  +2.23%  '
$ws.Range('D15').Value = '3.245.25'
$ws.Range('E15').Value = '  +0.83%  '
$ws.Range('D16').Value = '2.802.01'
$ws.Range('E16').Value = '  +0.89%  '
$ws.Range('E17').Value = '  +1.38%  '
$ws.Range('D18').Value = '51.777.47'
$ws.Range('E18').Value = '  +0.17%  '
$ws.Range('E19').Value = '  +1.29%  '
$ws.Range('D20').Value = "'3.19"
$ws.Range('E20').Value = '  +3.48%  '
$ws.Range('D21').Value = "'13.67"
$ws.Range('E21').Value = '  +3.39%  '
$ws.Range('E22').Value = '  +0.85%  '
$ws.Range('D23').Value = "'70.56"
$ws.Range('E23').Value = '  +0.45%  '
$ws.Range('D24').Value = "'268.75"
$ws.Range('E24').Value = '  +0.35%  '
$ws.Range('E25').Value = '  +1.27%  '
$ws.Range('E27').Value = '  -0.82%  '
$ws.Range('E28').Value = '  -3.40%  '
$ws.Range('D29').Value = "'38.89"
$ws.Range('E29').Value = '  +11.23%  '
$ws.Range('E30').Value = '  +1.82%  '
$ws.Range('E31').Value = '  +3.61%  '
$ws.Range('D32').Value = "'52.29"
$ws.Range('E32').Value = '  +0.71%  '
$ws.Range('E33').Value = '  -0.32%  '
$ws.Range('E34').Value = '  +8.68%  '
$ws.Range('E35').Value = '  +5.83%  '
$ws.Range('E36').Value = '  -0.57%  '
$ws.Range('E37').Value = '  +0.03%  '
$ws.Range('D38').Value = "'18.90"
$ws.Range('E38').Value = '  +0.05%  '
$ws.Range('E39').Value = '  +2.37%  '
$ws.Range('D40').Value = "'3.15"
$ws.Range('E40').Value = '  +0.57%  '
$ws.Range('E41').Value = '  +1.21%  '
$ws.Range('E42').Value = '  -1.81%  '
$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').Value = "'2.22"
$ws.Range('E43').Value = '  +0.31%  '
$ws.Range('B44').Value = 'Monero'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D44').Value = "'119.97"
$ws.Range('E44').Value = '  +0.18%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = "'22.01"
$ws.Range('E45').Value = '  +1.23%  '
$ws.Range('D46').Value = "'3.41"
$ws.Range('E46').Value = '  +4.28%  '
$ws.Range('D47').Value = '2.117.93'
$ws.Range('E47').Value = '  +1.44%  '
$ws.Range('D48').Value = "'2.43"
$ws.Range('E48').Value = '  +7.11%  '
$ws.Range('E49').Value = '  +0.27%  '
$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D50').Value = "'5.48"
$ws.Range('E50').Value = '  -1.39%  '
$ws.Range('B51').Value = 'TrustWalletToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D51').Value = "'1.37"
$ws.Range('E51').Value = '  +7.35%  '

foreach ($addr in $textForcedCells) {
    $ws.Range($addr).Style = "Normal"
}
